# Update cryptocurrency price/volume snapshot values (columns D and E)
# generated from commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.882.96"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.447.17"
$ws.Range("E3").Value = "  -2.80%  "
$cell = $ws.Range("D4")
$cell.Formula = "'0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.Formula = "'523.52"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$cell = $ws.Range("D6")
$cell.Formula = "'130.84"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "
$cell = $ws.Range("D7")
$cell.Formula = "'0.999"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$cell = $ws.Range("D8")
$cell.Formula = "'0.564"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "2.450.97"
$ws.Range("E9").Value = "  -2.62%  "
$cell = $ws.Range("D10")
$cell.Formula = "'0.0983"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  -2.09%  "
$cell = $ws.Range("D12")
$cell.Formula = "'4.94"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Formula = "'0.324"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "2.880.59"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "57.767.97"
$ws.Range("E15").Value = "  -1.05%  "
$cell = $ws.Range("D16")
$cell.Formula = "'21.71"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "2.443.60"
$ws.Range("E18").Value = "  -2.85%  "
$cell = $ws.Range("D19")
$cell.Formula = "'10.29"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -3.06%  "
$cell = $ws.Range("D20")
$cell.Formula = "'4.14"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "
$cell = $ws.Range("D21")
$cell.Formula = "'315.44"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "
$cell = $ws.Range("D22")
$cell.Formula = "'6.09"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.11%  "
$cell = $ws.Range("D23")
$cell.Formula = "'1.00"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$cell = $ws.Range("D24")
$cell.Formula = "'64.65"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  -1.76%  "
$cell = $ws.Range("D28")
$cell.Formula = "'7.23"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.82%  "
$cell = $ws.Range("D29")
$cell.Formula = "'174.42"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  -2.48%  "
$cell = $ws.Range("D31")
$cell.Formula = "'1.70"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.18%  "
$cell = $ws.Range("D32")
$cell.Formula = "'6.10"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "
$cell = $ws.Range("D33")
$cell.Formula = "'1.15"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("E34").Value = "  +0.01%  "
$cell = $ws.Range("D35")
$cell.Formula = "'0.997"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$cell = $ws.Range("D36")
$cell.Formula = "'17.82"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.38%  "
$cell = $ws.Range("D37")
$cell.Formula = "'1.19"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Formula = "'3.77"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -3.98%  "
$cell = $ws.Range("D39")
$cell.Formula = "'36.44"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$cell = $ws.Range("D40")
$cell.Formula = "'1.45"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.70%  "
$cell = $ws.Range("D41")
$cell.Formula = "'0.791"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Formula = "'3.42"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.73%  "
$cell = $ws.Range("D43")
$cell.Formula = "'263.72"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -5.05%  "
$cell = $ws.Range("D44")
$cell.Formula = "'0.586"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.73%  "
$cell = $ws.Range("D45")
$cell.Formula = "'4.81"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.40%  "
$cell = $ws.Range("D46")
$cell.Formula = "'0.0925"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.50%  "
$cell = $ws.Range("D47")
$cell.Formula = "'122.40"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -4.18%  "
$cell = $ws.Range("D48")
$cell.Formula = "'0.0494"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$cell = $ws.Range("D49")
$cell.Formula = "'0.0211"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$cell = $ws.Range("D50")
$cell.Formula = "'16.98"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.74%  "
$cell = $ws.Range("D51")
$cell.Formula = "'16.36"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.61%  "
